$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 24 - this pushes the existing "Vehicle" datatype
# block (previously rows 24-51) down to rows 25-52, matching the reference edit.
$ws.Rows("24").Insert()

# New field row for the Policy datatype: PaymentPlan[][] paymentMatrix
# (the string table records the value on first use, column C before B,
# to line up with the target shared-strings ordering)
$ws.Range("C23").Value = "paymentMatrix"
$ws.Range("B23").Value = "PaymentPlan[][]"

# New PaymentPlan datatype table
$ws.Range("B55").Value = "Datatype PaymentPlan"
$ws.Range("B56").Value = "PlanName"
$ws.Range("C56").Value = "name"

# New PlanName <String> enumeration datatype table
$ws.Range("B59").Value = "Datatype PlanName <String>"
$ws.Range("B60").Value = "ANNUAL"
$ws.Range("B61").Value = "NONANNUAL"

# Match the selection left behind by the edit (viewport scroll position is
# cosmetic UI state Excel recomputes from live scrolling and isn't otherwise
# settable here).
$ws.Range("G69").Select()
